# Update the cryptos worksheet with the latest scraped price / volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '92.120.43'
$ws.Range("E2").Value = '  +1.96%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.110.93'
$ws.Range("E3").Value = '  +0.75%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.25'
$ws.Range("E5").Value = '  -1.26%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '614.13'
$ws.Range("E6").Value = '  -0.80%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.09'
$ws.Range("E7").Value = '  -4.64%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.391'
$ws.Range("E8").Value = '  +7.68%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.105.37'
$ws.Range("E10").Value = '  +0.61%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.730'
$ws.Range("E11").Value = '  -0.97%  '

# Row 12
$ws.Range("E12").Value = '  -0.73%  '

# Row 13
$ws.Range("E13").Value = '  +2.22%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.51'
$ws.Range("E14").Value = '  +1.02%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.714.91'
$ws.Range("E15").Value = '  +1.62%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.21'
$ws.Range("E16").Value = '  -1.99%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.679.82'
$ws.Range("E17").Value = '  +0.29%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.104.06'
$ws.Range("E18").Value = '  +0.52%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.65'
$ws.Range("E19").Value = '  -2.36%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.79'
$ws.Range("E20").Value = '  +0.77%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.83'
$ws.Range("E21").Value = '  +0.48%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '447.60'
$ws.Range("E22").Value = '  +2.29%  '

# Row 23
$ws.Range("B23").Value = 'PEPE'
$ws.Range("C23").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000203'
$ws.Range("E23").Value = '  -2.47%  '

# Row 24
$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.30'
$ws.Range("E24").Value = '  +2.73%  '

# Row 25
$ws.Range("E25").Value = '  -0.69%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.03'
$ws.Range("E26").Value = '  -2.47%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.68'
$ws.Range("E27").Value = '  -1.09%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.265.60'
$ws.Range("E28").Value = '  +0.63%  '

# Row 29
$ws.Range("E29").Value = '  +0.11%  '

# Row 30
$ws.Range("E30").Value = '  +12.83%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.228'
$ws.Range("E31").Value = '  -7.25%  '

# Row 32
$ws.Range("E32").Value = '  -4.90%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.28'
$ws.Range("E33").Value = '  +0.99%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +55.34%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.167'
$ws.Range("E35").Value = '  -1.06%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.01'
$ws.Range("E36").Value = '  +0.91%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.19'
$ws.Range("E37").Value = '  +0.05%  '

# Row 38
$ws.Range("E38").Value = '  -6.27%  '

# Row 39
$ws.Range("E39").Value = '  +1.14%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '481.84'
$ws.Range("E40").Value = '  -2.07%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.31'
$ws.Range("E41").Value = '  +1.25%  '

# Row 42
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.47'
$ws.Range("E42").Value = '  -2.37%  '

# Row 43
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.434'
$ws.Range("E43").Value = '  +3.58%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.22'
$ws.Range("E44").Value = '  +0.24%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '159.00'
$ws.Range("E46").Value = '  +3.34%  '

# Row 47
$ws.Range("E47").Value = '  +0.34%  '

# Row 48
$ws.Range("E48").Value = '  +1.49%  '

# Row 49
$ws.Range("E49").Value = '  +2.56%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0328'
$ws.Range("E50").Value = '  +5.27%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.98'
$ws.Range("E51").Value = '  -0.61%  '
